$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "0.999") need to be
# forced to stay text-typed, matching the source inline-string cells, otherwise
# Excel auto-converts them to numeric cells (dropping formatting like trailing zeros).
$textGuardCells = @(
    "D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D20", "D21", "D22", "D23",
    "D24", "D27", "D28", "D31", "D34", "D35", "D36", "D37", "D38", "D41", "D43", "D44",
    "D46", "D48", "D50", "D51"
)
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data
$ws.Range("D2").Value = '68.619.40'
$ws.Range("E2").Value = '  +4.14%  '

$ws.Range("D3").Value = '3.370.26'

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '593.97'
$ws.Range("E5").Value = '  +6.41%  '

$ws.Range("D6").Value = '185.73'
$ws.Range("E6").Value = '  +0.41%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  +3.75%  '

$ws.Range("E9").Value = '  +4.04%  '

$ws.Range("D10").Value = '0.585'
$ws.Range("E10").Value = '  +1.52%  '

$ws.Range("D11").Value = '47.13'
$ws.Range("E11").Value = '  +2.97%  '

$ws.Range("D12").Value = '0.0000279'
$ws.Range("E12").Value = '  +6.91%  '

$ws.Range("D13").Value = '641.00'
$ws.Range("E13").Value = '  +12.78%  '

$ws.Range("D14").Value = '3.909.77'
$ws.Range("E14").Value = '  +1.52%  '

$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").Value = '68.761.65'
$ws.Range("E16").Value = '  +4.49%  '

$ws.Range("E17").Value = '  +1.85%  '

$ws.Range("D18").Value = '3.372.62'
$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("E19").Value = '  +1.34%  '

$ws.Range("D20").Value = '11.05'
$ws.Range("E20").Value = '  +2.12%  '

$ws.Range("D21").Value = '0.909'
$ws.Range("E21").Value = '  +2.23%  '

$ws.Range("D22").Value = '17.93'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").Value = '5.07'
$ws.Range("E23").Value = '  +1.97%  '

$ws.Range("D24").Value = '98.96'
$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("E25").Value = '  +3.81%  '

$ws.Range("E26").Value = '  +5.69%  '

$ws.Range("D27").Value = '9.80'
$ws.Range("E27").Value = '  +4.80%  '

$ws.Range("D28").Value = '32.95'
$ws.Range("E28").Value = '  +8.40%  '

$ws.Range("E29").Value = '  +2.30%  '

$ws.Range("E30").Value = '  +1.55%  '

$ws.Range("D31").Value = '611.16'
$ws.Range("E31").Value = '  +9.44%  '

$ws.Range("D32").Value = '3.986.58'
$ws.Range("E32").Value = '  +6.74%  '

$ws.Range("E33").Value = '  +0.15%  '

$ws.Range("D34").Value = '11.07'
$ws.Range("E34").Value = '  +2.41%  '

$ws.Range("D35").Value = '0.105'
$ws.Range("E35").Value = '  +2.42%  '

$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").Value = '56.09'
$ws.Range("E37").Value = '  +1.02%  '

$ws.Range("D38").Value = '2.76'
$ws.Range("E38").Value = '  +7.20%  '

$ws.Range("E39").Value = '  +6.27%  '

$ws.Range("E40").Value = '  +4.28%  '

$ws.Range("D41").Value = '33.52'
$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("D42").Value = '0.0₃0703'
$ws.Range("E42").Value = '  +2.97%  '

$ws.Range("D43").Value = '3.43'
$ws.Range("E43").Value = '  +3.06%  '

$ws.Range("D44").Value = '0.342'
$ws.Range("E44").Value = '  +3.17%  '

$ws.Range("E45").Value = '  +3.41%  '

$ws.Range("D46").Value = '0.129'
$ws.Range("E46").Value = '  +2.27%  '

$ws.Range("E47").Value = '  +3.27%  '

$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.47%  '

$ws.Range("E49").Value = '  +9.11%  '

$ws.Range("D50").Value = '131.17'
$ws.Range("E50").Value = '  +4.82%  '

$ws.Range("D51").Value = '7.79'
$ws.Range("E51").Value = '  +7.46%  '

# Restore the default (Normal) style on the guarded cells so only the value
# changed -- NumberFormat reverts to General/no explicit style, same as source.
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).Style = "Normal"
}
